# Applies updated NATMI ligand-receptor statistics ("Natmi following Dr Hou advice")
# to the Thbs1-Itga2b sheet. New values follow from:
#   - an updated expressing-cell count per cluster (both as sender/ligand side
#     and as target/receptor side)
#   - ligand statistics (avg/total expression & specificity) recomputed per sending cluster
#   - receptor statistics (avg/total expression & specificity) recomputed per target cluster
#   - edge statistics = ligand-side value * receptor-side value (average and total weights
#     and specificities)

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$clusters = @("ECs", "FAPs", "M1", "M2", "Neutro", "sCs")

# Updated "expressing cells" counts per cluster (used for both the ligand/E column
# when the cluster is the sender, and the receptor/K column when the cluster is the target)
$exprCount = @{
    "ECs"    = 2
    "FAPs"   = 3
    "M1"     = 3
    "M2"     = 3
    "Neutro" = 3
    "sCs"    = 2
}

# Updated ligand statistics per sending cluster
$ligandAvg = @{
    "ECs"    = 31.645482
    "FAPs"   = 83.22744366666667
    "M1"     = 524.5768889999999
    "M2"     = 181.4813383333333
    "Neutro" = 15.79677433333333
    "sCs"    = 53.1948375
}
$ligandAvgSpec = @{
    "ECs"    = 0.03555980726701226
    "FAPs"   = 0.09352209759714789
    "M1"     = 0.5894633891046084
    "M2"     = 0.2039293133121744
    "Neutro" = 0.01775072507139627
    "sCs"    = 0.05977466764766092
}
$ligandTotSpec = @{
    "ECs"    = 0.0244846141215985
    "FAPs"   = 0.09659160077758068
    "M1"     = 0.6088102578564109
    "M2"     = 0.210622508737405
    "Neutro" = 0.01833332435500452
    "sCs"    = 0.0411576941520005
}

# Updated receptor statistics per target cluster
$receptorAvg = @{
    "ECs"    = 1.2502825
    "FAPs"   = 2.596814
    "M1"     = 0.2737596666666667
    "M2"     = 0.7929023333333333
    "Neutro" = 0.879594
    "sCs"    = 1.8364275
}
$receptorAvgSpec = @{
    "ECs"    = 0.1638687485091313
    "FAPs"   = 0.3403524085884521
    "M1"     = 0.03588041420154535
    "M2"     = 0.103922044060685
    "Neutro" = 0.1152843201245645
    "sCs"    = 0.2406920645156217
}
$receptorTotSpec = @{
    "ECs"    = 0.1262743888059313
    "FAPs"   = 0.3934044114342388
    "M1"     = 0.04147322855600493
    "M2"     = 0.1201207617372036
    "Neutro" = 0.1332541182661091
    "sCs"    = 0.1854730912005122
}

$row = 2
foreach ($snd in $clusters) {
    $E = $exprCount[$snd]
    $G = $ligandAvg[$snd]
    $H = $G * $E
    $I = $ligandAvgSpec[$snd]
    $J = $ligandTotSpec[$snd]

    foreach ($tgt in $clusters) {
        $K = $exprCount[$tgt]
        $M = $receptorAvg[$tgt]
        $N = $M * $K
        $O = $receptorAvgSpec[$tgt]
        $P = $receptorTotSpec[$tgt]

        $Q = $G * $M
        $R = $H * $N
        $S = $I * $O
        $T = $J * $P

        $ws.Cells.Item($row, 5).Value  = $E
        $ws.Cells.Item($row, 7).Value  = $G
        $ws.Cells.Item($row, 8).Value  = $H
        $ws.Cells.Item($row, 9).Value  = $I
        $ws.Cells.Item($row, 10).Value = $J
        $ws.Cells.Item($row, 11).Value = $K
        $ws.Cells.Item($row, 13).Value = $M
        $ws.Cells.Item($row, 14).Value = $N
        $ws.Cells.Item($row, 15).Value = $O
        $ws.Cells.Item($row, 16).Value = $P
        $ws.Cells.Item($row, 17).Value = $Q
        $ws.Cells.Item($row, 18).Value = $R
        $ws.Cells.Item($row, 19).Value = $S
        $ws.Cells.Item($row, 20).Value = $T

        $row = $row + 1
    }
}

Write-Host "Updated rows 2..$($row - 1)"
